$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update taxon labels to the corrected casing
$ws.Range("A9").Value = "Dragonflies & damselflies"
$ws.Range("A10").Value = "Dragonflies & damselflies"

$ws.Range("A11").Value = "Freshwater crabs"
$ws.Range("A12").Value = "Freshwater crabs"

$ws.Range("A13").Value = "Freshwater fishes"
$ws.Range("A14").Value = "Freshwater fishes"

# Reset the view: scroll back to the top and select A14
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A14").Select()
